$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("D2").Value = "India"
$ws.Range("D3").Value = "Russia"
$ws.Range("D4").Value = "China"
